$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.467854666666666
$ws.Range("H2").Value = 25.403564
$ws.Range("I2").Value = 0.1523462492674013
$ws.Range("J2").Value = 0.1523462492674013
$ws.Range("M2").Value = 8.533046666666666
$ws.Range("N2").Value = 25.59914
$ws.Range("O2").Value = 0.2932132236642383
$ws.Range("P2").Value = 0.2932132236642383
$ws.Range("Q2").Value = 72.25659903721777
$ws.Range("R2").Value = 650.30939133496
$ws.Range("S2").Value = 0.04466993486085033
$ws.Range("T2").Value = 0.04466993486085034

# Row 3
$ws.Range("G3").Value = 8.467854666666666
$ws.Range("H3").Value = 25.403564
$ws.Range("I3").Value = 0.1523462492674013
$ws.Range("J3").Value = 0.1523462492674013
$ws.Range("O3").Value = 0.3119288965200195
$ws.Range("P3").Value = 0.3119288965200194
$ws.Range("Q3").Value = 76.8687064052009
$ws.Range("R3").Value = 691.8183576468081
$ws.Range("S3").Value = 0.04752119742294431
$ws.Range("T3").Value = 0.04752119742294431

# Row 4
$ws.Range("G4").Value = 8.467854666666666
$ws.Range("H4").Value = 25.403564
$ws.Range("I4").Value = 0.1523462492674013
$ws.Range("J4").Value = 0.1523462492674013
$ws.Range("O4").Value = 0.3948578798157423
$ws.Range("P4").Value = 0.3948578798157423
$ws.Range("Q4").Value = 97.30491395300531
$ws.Range("R4").Value = 875.7442255770479
$ws.Range("S4").Value = 0.06015511698360665
$ws.Range("T4").Value = 0.06015511698360666

# Row 5
$ws.Range("I5").Value = 0.1669927598427297
$ws.Range("J5").Value = 0.1669927598427297
$ws.Range("M5").Value = 8.533046666666666
$ws.Range("N5").Value = 25.59914
$ws.Range("O5").Value = 0.2932132236642383
$ws.Range("P5").Value = 0.2932132236642383
$ws.Range("Q5").Value = 79.20332104071333
$ws.Range("R5").Value = 712.82988936642
$ws.Range("S5").Value = 0.04896448544207475
$ws.Range("T5").Value = 0.04896448544207475

# Row 6
$ws.Range("I6").Value = 0.1669927598427297
$ws.Range("J6").Value = 0.1669927598427297
$ws.Range("O6").Value = 0.3119288965200195
$ws.Range("P6").Value = 0.3119288965200194
$ws.Range("S6").Value = 0.05208986730457531
$ws.Range("T6").Value = 0.0520898673045753

# Row 7
$ws.Range("I7").Value = 0.1669927598427297
$ws.Range("J7").Value = 0.1669927598427297
$ws.Range("O7").Value = 0.3948578798157423
$ws.Range("P7").Value = 0.3948578798157423
$ws.Range("S7").Value = 0.06593840709607969
$ws.Range("T7").Value = 0.06593840709607969

# Row 8
$ws.Range("I8").Value = 0.680660990889869
$ws.Range("J8").Value = 0.680660990889869
$ws.Range("M8").Value = 8.533046666666666
$ws.Range("N8").Value = 25.59914
$ws.Range("O8").Value = 0.2932132236642383
$ws.Range("P8").Value = 0.2932132236642383
$ws.Range("Q8").Value = 322.8320259639533
$ws.Range("R8").Value = 2905.48823367558
$ws.Range("S8").Value = 0.1995788033613132
$ws.Range("T8").Value = 0.1995788033613132

# Row 9
$ws.Range("I9").Value = 0.680660990889869
$ws.Range("J9").Value = 0.680660990889869
$ws.Range("O9").Value = 0.3119288965200195
$ws.Range("P9").Value = 0.3119288965200194
$ws.Range("S9").Value = 0.2123178317924999
$ws.Range("T9").Value = 0.2123178317924998

# Row 10
$ws.Range("I10").Value = 0.680660990889869
$ws.Range("J10").Value = 0.680660990889869
$ws.Range("O10").Value = 0.3948578798157423
$ws.Range("P10").Value = 0.3948578798157423
$ws.Range("S10").Value = 0.2687643557360559
$ws.Range("T10").Value = 0.2687643557360559
